$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-5 (previously held by rows 6-9 "Compass Error" content)
$ws.Range("B2").Value = "Compass Error Compass Error Compass disconnected ."
$ws.Range("C2").Value = "Compass Error"
$ws.Range("D2").Value = "0-1"

$ws.Range("B3").Value = "Compass Error Compass Error Compass disconnected ."
$ws.Range("C3").Value = "Compass Error"
$ws.Range("D3").Value = "2-3"

$ws.Range("B4").Value = "Compass Error Compass Error Compass disconnected ."
$ws.Range("C4").Value = "Compass disconnected"
$ws.Range("D4").Value = "4-5"

$ws.Range("B5").Value = "Compass Error Compass Error Compass disconnected ."
$ws.Range("C5").Value = "Compass Error Compass Error Compass disconnected"
$ws.Range("D5").Value = "0-5"

# New values for rows 6-9 (previously held by rows 2-5 "Critical low battery" content)
$ws.Range("B6").Value = "Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."
$ws.Range("C6").Value = "Critical low battery"
$ws.Range("D6").Value = "0-2"

$ws.Range("B7").Value = "Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."
$ws.Range("C7").Value = "Aircraft in Auto Power Off Protection"
$ws.Range("D7").Value = "3-8"

$ws.Range("B8").Value = "Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."
$ws.Range("C8").Value = "Forced landing in progress"
$ws.Range("D8").Value = "9-12"

$ws.Range("B9").Value = "Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."
$ws.Range("C9").Value = "Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress"
$ws.Range("D9").Value = "0-12"
